$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2006
$ws.Range("J7").Value = 2006
$ws.Range("L7").Value = 2006
$ws.Range("N7").Value = -2230
$ws.Range("H14").Value = 2006
$ws.Range("J14").Value = 2006
$ws.Range("L14").Value = 2006
$ws.Range("N14").Value = -2388
$ws.Range("H18").Value = 959.3077
$ws.Range("I18").Value = 959.3077
$ws.Range("K18").Value = 959.3077
$ws.Range("M18").Value = -675.3077
$ws.Range("H28").Value = 263.15152
$ws.Range("I28").Value = 252
$ws.Range("J28").Value = 304.57144
$ws.Range("K28").Value = 252
$ws.Range("L28").Value = 304.57144
$ws.Range("M28").Value = 233
$ws.Range("N28").Value = -1274.57144
$ws.Range("H44").Value = 19500
$ws.Range("J44").Value = 19500
$ws.Range("L44").Value = 19500
$ws.Range("N44").Value = -20424
$ws.Range("H51").Value = 2216.6667
$ws.Range("I51").Value = 1625
$ws.Range("J51").Value = 3400
$ws.Range("K51").Value = 1625
$ws.Range("L51").Value = 3400
$ws.Range("M51").Value = -1141
$ws.Range("N51").Value = -4368
$ws.Range("H92").Value = 404.44446
$ws.Range("I92").Value = 205.71428
$ws.Range("J92").Value = 1100
$ws.Range("K92").Value = 205.71428
$ws.Range("L92").Value = 1100
$ws.Range("M92").Value = 1042.28572
$ws.Range("N92").Value = -3596
$ws.Range("H94").Value = 6725.625
$ws.Range("I94").Value = 6725.625
$ws.Range("K94").Value = 6725.625
$ws.Range("M94").Value = -6274.625
$ws.Range("H96").Value = 507.33334
$ws.Range("I96").Value = 368.25
$ws.Range("J96").Value = 618.6
$ws.Range("K96").Value = 1104.75
$ws.Range("L96").Value = 1855.8
$ws.Range("M96").Value = 268.25
$ws.Range("N96").Value = -4601.8
$ws.Range("H100").Value = 3483.75
$ws.Range("J100").Value = 3580.5
$ws.Range("L100").Value = 3580.5
$ws.Range("N100").Value = -4662.5
$ws.Range("H129").Value = 1004.5
$ws.Range("I129").Value = 816.8333
$ws.Range("J129").Value = 1028.9783
$ws.Range("K129").Value = 2450.4999
$ws.Range("L129").Value = 3086.9349
$ws.Range("M129").Value = 2549.5001
$ws.Range("N129").Value = -13086.9349
$ws.Range("H132").Value = 2916.9143
$ws.Range("I132").Value = 2351.0476
$ws.Range("K132").Value = 7053.1428
$ws.Range("M132").Value = -4523.1428
$ws.Range("H138").Value = 3028.3718
$ws.Range("I138").Value = 1289.1428
$ws.Range("J138").Value = 3669.1404
$ws.Range("K138").Value = 3867.4284
$ws.Range("L138").Value = 11007.4212
$ws.Range("M138").Value = 1272.5716
$ws.Range("N138").Value = -21287.4212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H19").Value = 10502.667
$ws.Range("I19").Value = 754
$ws.Range("K19").Value = 754
$ws.Range("M19").Value = -525
$ws.Range("H102").Value = 1646.7222
$ws.Range("I102").Value = 1510
$ws.Range("J102").Value = 2002.2
$ws.Range("K102").Value = 1510
$ws.Range("L102").Value = 2002.2
$ws.Range("M102").Value = 112
$ws.Range("N102").Value = -5246.2
$ws.Range("H122").Value = 1397.3823
$ws.Range("I122").Value = 1096.5769
$ws.Range("J122").Value = 2375
$ws.Range("K122").Value = 3289.7307
$ws.Range("L122").Value = 7125
$ws.Range("M122").Value = -839.7307000000001
$ws.Range("N122").Value = -12025

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3751.8333
$ws.Range("I10").Value = 502.2
$ws.Range("K10").Value = 502.2
$ws.Range("M10").Value = -363.2
$ws.Range("H58").Value = 2005.3636
$ws.Range("I58").Value = 2107.9
$ws.Range("K58").Value = 2107.9
$ws.Range("M58").Value = -1904.9
$ws.Range("H94").Value = 20497.125
$ws.Range("I94").Value = 1660
$ws.Range("J94").Value = 23188.143
$ws.Range("K94").Value = 1660
$ws.Range("L94").Value = 23188.143
$ws.Range("M94").Value = -1209
$ws.Range("N94").Value = -24090.143
$ws.Range("H99").Value = 3053.4119
$ws.Range("I99").Value = 3166.1667
$ws.Range("J99").Value = 2782.8
$ws.Range("K99").Value = 3166.1667
$ws.Range("L99").Value = 2782.8
$ws.Range("M99").Value = -1668.1667
$ws.Range("N99").Value = -5778.8
$ws.Range("H126").Value = 3053.4119
$ws.Range("I126").Value = 3166.1667
$ws.Range("J126").Value = 2782.8
$ws.Range("K126").Value = 9498.500100000001
$ws.Range("L126").Value = 8348.400000000001
$ws.Range("M126").Value = -7028.500100000001
$ws.Range("N126").Value = -13288.4
$ws.Range("H136").Value = 2005.3636
$ws.Range("I136").Value = 2107.9
$ws.Range("K136").Value = 6323.700000000001
$ws.Range("M136").Value = -3773.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1530.5652
$ws.Range("J5").Value = 1291.6666
$ws.Range("L5").Value = 3874.9998
$ws.Range("N5").Value = -4098.9998
$ws.Range("H115").Value = 1816.8182
$ws.Range("I115").Value = 372.2
$ws.Range("J115").Value = 3020.6667
$ws.Range("K115").Value = 1116.6
$ws.Range("L115").Value = 9062.000100000001
$ws.Range("M115").Value = 58.40000000000009
$ws.Range("N115").Value = -11412.0001
$ws.Range("H131").Value = 14707314
$ws.Range("J131").Value = 15152977
$ws.Range("L131").Value = 45458931
$ws.Range("N131").Value = -45469011
$ws.Range("H135").Value = 1530.5652
$ws.Range("J135").Value = 1291.6666
$ws.Range("L135").Value = 11624.9994
$ws.Range("N135").Value = -16694.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 50000
$ws.Range("I29").Value = 50000
$ws.Range("K29").Value = 50000
$ws.Range("M29").Value = -49710
$ws.Range("H122").Value = 1794.4117
$ws.Range("I122").Value = 1269.6154
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 3808.8462
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -1358.8462
$ws.Range("N122").Value = -15400
$ws.Range("H123").Value = 14123.375
$ws.Range("J123").Value = 14123.375
$ws.Range("L123").Value = 14123.375
$ws.Range("N123").Value = -19023.375
$ws.Range("H125").Value = 48331.5
$ws.Range("J125").Value = 48331.5
$ws.Range("L125").Value = 48331.5
$ws.Range("N125").Value = -53251.5
$ws.Range("H126").Value = 3359.8
$ws.Range("I126").Value = 2899.5
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 8698.5
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -6228.5
$ws.Range("N126").Value = -15940.0001
$ws.Range("H131").Value = 44243.5
$ws.Range("J131").Value = 44243.5
$ws.Range("L131").Value = 44243.5
$ws.Range("N131").Value = -54323.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1538
$ws.Range("I16").Value = 1297.5
$ws.Range("K16").Value = 1297.5
$ws.Range("M16").Value = -1127.5
$ws.Range("H43").Value = 35000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20386
$ws.Range("H93").Value = 1285.7142
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 1666.6666
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 1666.6666
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -4162.6666
$ws.Range("H100").Value = 4026.2173
$ws.Range("I100").Value = 4535.4707
$ws.Range("J100").Value = 2583.3333
$ws.Range("K100").Value = 4535.4707
$ws.Range("L100").Value = 2583.3333
$ws.Range("M100").Value = -3994.4707
$ws.Range("N100").Value = -3665.3333
$ws.Range("H122").Value = 18754468
$ws.Range("I122").Value = 15628863
$ws.Range("J122").Value = 25005674
$ws.Range("K122").Value = 46886589
$ws.Range("L122").Value = 75017022
$ws.Range("M122").Value = -46884139
$ws.Range("N122").Value = -75021922

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 28900
$ws.Range("I29").Value = 28900
$ws.Range("K29").Value = 28900
$ws.Range("M29").Value = -28610
$ws.Range("H96").Value = 1475
$ws.Range("I96").Value = 1440
$ws.Range("J96").Value = 1533.3334
$ws.Range("K96").Value = 1440
$ws.Range("L96").Value = 1533.3334
$ws.Range("M96").Value = -67
$ws.Range("N96").Value = -4279.3334
$ws.Range("H114").Value = 80600
$ws.Range("J114").Value = 80600
$ws.Range("L114").Value = 80600
$ws.Range("N114").Value = -89278
$ws.Range("H115").Value = 53266.668
$ws.Range("J115").Value = 53266.668
$ws.Range("L115").Value = 53266.668
$ws.Range("N115").Value = -56400.668
